# 12.05 PSP(김지환, 함형우) 수정, SRS 수정, TestCase 수정
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 김지환 sheet: fill in row 30 (new PSP time-log entry)
# ---------------------------------------------------------------
$wsJH = $wb.Worksheets.Item("김지환")

$wsJH.Range("A30").Value = "19.12.04"
$wsJH.Range("B30").Value = 0.95833333333333337
$wsJH.Range("C30").Value = 0.1111111111111111
$wsJH.Range("D30").Value = 0
$wsJH.Range("E30").Value = 220
$wsJH.Range("F30").Value = "DB boundary 함수 구현, 피드백 부분 Test Scenario 작성, SRS 수정 "

$wsJH.Rows.Item(30).RowHeight = 15

# ---------------------------------------------------------------
# 함형우 sheet: fill in rows 13-18 (new PSP time-log entries)
# ---------------------------------------------------------------
$wsHW = $wb.Worksheets.Item("함형우")

$wsHW.Range("A13").Value = "19.11.13"
$wsHW.Range("B13").Value = 0.41666666666666669
$wsHW.Range("C13").Value = 0.58333333333333337
$wsHW.Range("D13").Value = 0
$wsHW.Range("E13").Value = 120
$wsHW.Range("F13").Value = "모임등록 UseCase 코드작업"

$wsHW.Range("A14").Value = "19.11.15"
$wsHW.Range("B14").Value = 0.75
$wsHW.Range("C14").Value = 0.84722222222222221
$wsHW.Range("D14").Value = 0
$wsHW.Range("E14").Value = 140
$wsHW.Range("F14").Value = "모임등록 UseCase 코드작업"

$wsHW.Range("A15").Value = "19.11.21"
$wsHW.Range("B15").Value = 0.83333333333333337
$wsHW.Range("C15").Value = 0.95138888888888884
$wsHW.Range("D15").Value = 0
$wsHW.Range("E15").Value = 170
$wsHW.Range("F15").Value = "Everytime 연동 관련 코드 작업"

$wsHW.Range("A16").Value = "19.11.23"
$wsHW.Range("B16").Value = 0.625
$wsHW.Range("C16").Value = 0.65972222222222221
$wsHW.Range("D16").Value = 0
$wsHW.Range("E16").Value = 50
$wsHW.Range("F16").Value = "모임등록 UseCase 코드작업"

$wsHW.Range("A17").Value = "19.12.1"
$wsHW.Range("B17").Value = 0.5
$wsHW.Range("C17").Value = 0.58333333333333337
$wsHW.Range("D17").Value = 0
$wsHW.Range("E17").Value = 120
$wsHW.Range("F17").Value = "Everytime 연동 관련 코드 작업"

$wsHW.Range("A18").Value = "19.12.3"
$wsHW.Range("B18").Value = 0.78472222222222221
$wsHW.Range("C18").Value = 0.86111111111111116
$wsHW.Range("D18").Value = 0
$wsHW.Range("E18").Value = 110
$wsHW.Range("F18").Value = "코드 수정"

# ---------------------------------------------------------------
# Selection / active-sheet bookkeeping: the edit moved the user's
# working selection from 김지환!F29 to 김지환!F30, then on to
# 함형우!F37, leaving 함형우 as the active (tabSelected) sheet.
# ---------------------------------------------------------------
$wsJH.Activate()
$wsJH.Range("F30").Select()

$wsHW.Activate()
$wsHW.Range("F37").Select()
